# Auto-generated edit script: updates crypto price/volume table cells
# to match the latest GitHub Actions scrape snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $value) {
    $cell = $ws.Range($ref)
    # Force text storage so numeric-looking strings (prices like
    # "576.36" or "0.290") are kept verbatim instead of being
    # coerced into floating point numbers by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    # Reset the cell style back to Normal/General so we do not
    # leave a stray "Text" number-format style behind.
    $cell.Style = "Normal"
}

Set-TextCell 'D2' '65.176.76'
Set-TextCell 'E2' '  +1.76%  '
Set-TextCell 'D3' '3.185.55'
Set-TextCell 'E3' '  +4.11%  '
Set-TextCell 'E4' '  +0.00%  '
Set-TextCell 'D5' '576.36'
Set-TextCell 'E5' '  +2.90%  '
Set-TextCell 'D6' '151.41'
Set-TextCell 'E6' '  +5.56%  '
Set-TextCell 'E7' '  +0.05%  '
Set-TextCell 'D8' '3.184.07'
Set-TextCell 'E8' '  +4.11%  '
Set-TextCell 'D9' '0.529'
Set-TextCell 'E9' '  +3.47%  '
Set-TextCell 'D10' '0.164'
Set-TextCell 'E10' '  +5.21%  '
Set-TextCell 'D11' '6.24'
Set-TextCell 'E11' '  +2.47%  '
Set-TextCell 'D12' '0.509'
Set-TextCell 'E12' '  +4.99%  '
Set-TextCell 'D13' '0.0000277'
Set-TextCell 'E13' '  +19.88%  '
Set-TextCell 'D14' '38.27'
Set-TextCell 'E14' '  +8.13%  '
Set-TextCell 'D15' '3.704.07'
Set-TextCell 'E15' '  +4.02%  '
Set-TextCell 'D16' '65.283.18'
Set-TextCell 'E16' '  +1.86%  '
Set-TextCell 'D17' '3.184.34'
Set-TextCell 'E17' '  +3.88%  '
Set-TextCell 'D18' '7.23'
Set-TextCell 'E18' '  +7.14%  '
Set-TextCell 'D20' '515.10'
Set-TextCell 'E20' '  +7.83%  '
Set-TextCell 'D21' '14.98'
Set-TextCell 'E21' '  +7.25%  '
Set-TextCell 'D22' '0.737'
Set-TextCell 'E22' '  +8.35%  '
Set-TextCell 'D23' '15.59'
Set-TextCell 'E23' '  +9.13%  '
Set-TextCell 'D24' '7.87'
Set-TextCell 'E24' '  +4.17%  '
Set-TextCell 'D25' '85.20'
Set-TextCell 'E25' '  +3.72%  '
Set-TextCell 'E26' '  +0.16%  '
Set-TextCell 'D27' '9.10'
Set-TextCell 'E27' '  +13.54%  '
Set-TextCell 'E28' '  +4.41%  '
Set-TextCell 'E29' '  +8.48%  '
Set-TextCell 'D30' '28.20'
Set-TextCell 'E30' '  +7.17%  '
Set-TextCell 'D31' '2.79'
Set-TextCell 'E31' '  +14.56%  '
Set-TextCell 'E32' '  +7.56%  '
Set-TextCell 'E33' '  -0.05%  '
Set-TextCell 'E34' '  +11.51%  '
Set-TextCell 'D35' '6.76'
Set-TextCell 'E35' '  +8.25%  '
Set-TextCell 'D36' '55.76'
Set-TextCell 'E36' '  +1.52%  '
Set-TextCell 'D37' '0.0897'
Set-TextCell 'E37' '  +10.59%  '
Set-TextCell 'D38' '480.77'
Set-TextCell 'E38' '  +8.33%  '
Set-TextCell 'D39' '3.12'
Set-TextCell 'E39' '  +9.99%  '
Set-TextCell 'D40' '0.0424'
Set-TextCell 'E40' '  +3.85%  '
Set-TextCell 'D41' '3.145.40'
Set-TextCell 'E41' '  +5.08%  '
Set-TextCell 'D42' '8.67'
Set-TextCell 'E42' '  +5.08%  '
Set-TextCell 'E43' '  +4.77%  '
Set-TextCell 'B44' 'TheGraph'
Set-TextCell 'C44' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell 'D44' '0.290'
Set-TextCell 'E44' '  +10.91%  '
Set-TextCell 'B45' 'Fetch.AI'
Set-TextCell 'C45' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell 'D45' '2.48'
Set-TextCell 'E45' '  +15.06%  '
Set-TextCell 'D46' '29.53'
Set-TextCell 'E46' '  +6.29%  '
Set-TextCell 'D47' '0.0₃0613'
Set-TextCell 'E47' '  +17.98%  '
Set-TextCell 'E48' '  -0.07%  '
Set-TextCell 'E49' '  +2.37%  '
Set-TextCell 'E50' '  +11.31%  '
Set-TextCell 'D51' '122.18'
Set-TextCell 'E51' '  +3.12%  '
